$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old content in case old range was larger than new header area (not needed here, but safe)
# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows: Indice, Distancia, max, min, Tempo
$data = @(
    @(0, 1057, 1057, 1057, 0.01037770907084147),
    @(1, 1040, 1040, 1040, 0.01055965423583984),
    @(2, 973, 973, 973, 0.01231667200724284),
    @(3, 1224, 1224, 1224, 0.01220448017120361),
    @(4, 883, 883, 883, 0.01216433842976888),
    @(5, 1040, 1040, 1040, 0.0125093142191569),
    @(6, 1053, 1053, 1053, 0.01237522761027018),
    @(7, 957, 957, 957, 0.01046140193939209),
    @(8, 886, 886, 886, 0.01229832967122396),
    @(9, 1049, 1049, 1049, 0.0119659423828125)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
